# Add 2022-Q4 data:
#  - New sheet "2022-Q4" inserted right after "总计" (before "2022-Q3"),
#    built from a copy of "2022-Q3" so header/column-A styling matches the
#    other quarter sheets, then its content is replaced with the new data.
#  - "总计" summary sheet gets a new row inserted for 2022-Q4 (at the top
#    of the data, pushing the older quarters down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q4" worksheet from a copy of "2022-Q3" (keeps the
#    same header/border/bold styling used by every quarter sheet).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item("总计").Index + 1
$q4 = $wb.Worksheets.Item($q4)
$q4.Name = "2022-Q4"

# Remove the old Q3 rows (6-11) that got copied along with the sheet -
# 2022-Q4 only has 4 holdings.
$q4.Rows("6:11").Delete()

# Text-looking numeric columns (B..G) must stay plain text, like every
# other quarter sheet (only column H - rank - is a real number).
$q4.Range("B2:G5").NumberFormat = "@"

$q4.Cells.Item(2,2).Value = "012526"
$q4.Cells.Item(2,3).Value = "广发盛锦混合A"
$q4.Cells.Item(2,4).Value = "24.00"
$q4.Cells.Item(2,5).Value = "93.13"
$q4.Cells.Item(2,6).Value = "3.31"
$q4.Cells.Item(2,7).Value = "0.7944"
$q4.Cells.Item(2,8).Value = 10

$q4.Cells.Item(3,2).Value = "012527"
$q4.Cells.Item(3,3).Value = "广发盛锦混合C"
$q4.Cells.Item(3,4).Value = "1.14"
$q4.Cells.Item(3,5).Value = "93.13"
$q4.Cells.Item(3,6).Value = "3.31"
$q4.Cells.Item(3,7).Value = "0.0377"
$q4.Cells.Item(3,8).Value = 10

$q4.Cells.Item(4,2).Value = "015921"
$q4.Cells.Item(4,3).Value = "申万菱信国证2000指数增强A"
$q4.Cells.Item(4,4).Value = "0.21"
$q4.Cells.Item(4,5).Value = "94.00"
$q4.Cells.Item(4,6).Value = "0.51"
$q4.Cells.Item(4,7).Value = "0.0011"
$q4.Cells.Item(4,8).Value = 8

$q4.Cells.Item(5,2).Value = "015922"
$q4.Cells.Item(5,3).Value = "申万菱信国证2000指数增强C"
$q4.Cells.Item(5,4).Value = "0.08"
$q4.Cells.Item(5,5).Value = "94.00"
$q4.Cells.Item(5,6).Value = "0.51"
$q4.Cells.Item(5,7).Value = "0.0004"
$q4.Cells.Item(5,8).Value = 8

# ---------------------------------------------------------------------
# 2) Update "总计": insert a new top data row for 2022-Q4, pushing the
#    existing quarters down (2021-Q4 now lands on row 6).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()
$total.Range("A2:D2").ClearFormats()

# Column A carries the same bold/bordered style as the other index cells
# (A3:A6) - copy it across rather than guessing at the underlying xf.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 0.83

# Renumber the index column (A) for the rows that shifted down.
for ($r = 3; $r -le 6; $r++) {
    $total.Cells.Item($r,1).Value = $r - 2
}
